$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column E (dni_ciu) for rows 2-15, and corresponding
# recalculated column G (PORC_AVANCE = E / D * 100).
$updates = @(
    @{ Row = 2;  E = 1525; G = 26.11301369863014 },
    @{ Row = 3;  E = 1105; G = 36.83333333333334 },
    @{ Row = 4;  E = 1199; G = 27.18820861678005 },
    @{ Row = 5;  E = 1371; G = 27.42 },
    @{ Row = 6;  E = 1164; G = 14.55 },
    @{ Row = 7;  E = 1327; G = 22.51442144553783 },
    @{ Row = 8;  E = 1195; G = 17.07142857142857 },
    @{ Row = 9;  E = 1495; G = 37.375 },
    @{ Row = 10; E = 1255; G = 21.38718473074302 },
    @{ Row = 11; E = 2400; G = 43.05705059203444 },
    @{ Row = 12; E = 963;  G = 16.05 },
    @{ Row = 13; E = 1425; G = 28.5 },
    @{ Row = 14; E = 2497; G = 40.65450993161836 },
    @{ Row = 15; E = 1972; G = 39.24378109452736 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 7).Value = $u.G
}
